# Apply the daily TGP price-table roll-forward update.
# For each affected row: column A gets the new Effective Date serial,
# and columns D-G get the refreshed price figures (only the columns
# that actually carry a value in that row are touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 8; Cells = @{ "A" = 45954; "D" = 159.86000000000001; "E" = 157.69999999999999; "F" = 167.7; "G" = 157.86000000000001 } }
    @{ Row = 9; Cells = @{ "A" = 45954; "D" = 159.86000000000001; "E" = 157.69999999999999; "F" = 167.7; "G" = 157.86000000000001 } }
    @{ Row = 10; Cells = @{ "A" = 45954; "D" = 162.16; "E" = 159.91; "F" = 169.91; "G" = 160.38 } }
    @{ Row = 11; Cells = @{ "A" = 45953; "D" = 159.62; "E" = 157.68; "F" = 167.68; "G" = 157.83000000000001 } }
    @{ Row = 12; Cells = @{ "A" = 45953; "D" = 159.62; "E" = 157.68; "F" = 167.68; "G" = 157.83000000000001 } }
    @{ Row = 13; Cells = @{ "A" = 45953; "D" = 161.88999999999999; "E" = 159.75; "F" = 169.75; "G" = 160.22 } }
    @{ Row = 17; Cells = @{ "A" = 45954; "D" = 165.59; "E" = 162.78; "F" = 172.78 } }
    @{ Row = 18; Cells = @{ "A" = 45953; "D" = 165.32; "E" = 162.80000000000001; "F" = 172.8 } }
    @{ Row = 22; Cells = @{ "A" = 45954; "D" = 160.88999999999999; "E" = 158.81; "F" = 168.41; "G" = 159.99 } }
    @{ Row = 23; Cells = @{ "A" = 45954; "D" = 166.93; "E" = 163.63; "F" = 173.63 } }
    @{ Row = 24; Cells = @{ "A" = 45954; "D" = 166.73; "E" = 163.84; "F" = 173.84 } }
    @{ Row = 25; Cells = @{ "A" = 45954; "D" = 167.56; "E" = 163.22999999999999; "F" = 173.23; "G" = 163.05000000000001 } }
    @{ Row = 26; Cells = @{ "A" = 45954; "D" = 166.28; "E" = 164.76; "F" = 174.76 } }
    @{ Row = 27; Cells = @{ "A" = 45953; "D" = 160.54; "E" = 158.75; "F" = 168.35; "G" = 159.93 } }
    @{ Row = 28; Cells = @{ "A" = 45953; "D" = 166.65; "E" = 163.57; "F" = 173.57 } }
    @{ Row = 29; Cells = @{ "A" = 45953; "D" = 166.46; "E" = 163.78; "F" = 173.78 } }
    @{ Row = 30; Cells = @{ "A" = 45953; "D" = 167.29; "E" = 163.16999999999999; "F" = 173.17; "G" = 163 } }
    @{ Row = 31; Cells = @{ "A" = 45953; "D" = 166.01; "E" = 164.71; "F" = 174.71 } }
    @{ Row = 35; Cells = @{ "A" = 45954; "D" = 160.4; "E" = 157.13; "F" = 166.13 } }
    @{ Row = 36; Cells = @{ "A" = 45953; "D" = 160.13; "E" = 157.07; "F" = 166.07 } }
    @{ Row = 40; Cells = @{ "A" = 45954; "D" = 166.05; "E" = 162.51; "F" = 172.51 } }
    @{ Row = 41; Cells = @{ "A" = 45954; "D" = 165.77; "E" = 162.93; "F" = 172.93 } }
    @{ Row = 42; Cells = @{ "A" = 45953; "D" = 165.78; "E" = 162.54; "F" = 172.54 } }
    @{ Row = 43; Cells = @{ "A" = 45953; "D" = 165.5; "E" = 162.96; "F" = 172.96 } }
    @{ Row = 47; Cells = @{ "A" = 45954; "D" = 159.88; "E" = 158.61000000000001; "F" = 168.61 } }
    @{ Row = 48; Cells = @{ "A" = 45954; "D" = 159.86000000000001; "E" = 158.78; "F" = 168.78 } }
    @{ Row = 49; Cells = @{ "A" = 45953; "D" = 160.19999999999999; "E" = 158.82; "F" = 168.82 } }
    @{ Row = 50; Cells = @{ "A" = 45953; "D" = 160.18; "E" = 158.99; "F" = 168.99 } }
    @{ Row = 54; Cells = @{ "A" = 45954; "D" = 176.24; "E" = 172.88; "F" = 182.88 } }
    @{ Row = 55; Cells = @{ "A" = 45954; "D" = 163.9; "E" = 170.3; "F" = 180.3 } }
    @{ Row = 56; Cells = @{ "A" = 45954; "D" = 166.18 } }
    @{ Row = 57; Cells = @{ "A" = 45954; "D" = 165.85; "E" = 164.57 } }
    @{ Row = 58; Cells = @{ "A" = 45954; "D" = 161.76; "E" = 160.62; "F" = 170.62 } }
    @{ Row = 59; Cells = @{ "A" = 45954; "D" = 168.58; "E" = 171.08 } }
    @{ Row = 60; Cells = @{ "A" = 45953; "D" = 175.96; "E" = 172.84; "F" = 182.84 } }
    @{ Row = 61; Cells = @{ "A" = 45953; "D" = 163.63; "E" = 170.13; "F" = 180.13 } }
    @{ Row = 62; Cells = @{ "A" = 45953; "D" = 165.91 } }
    @{ Row = 63; Cells = @{ "A" = 45953; "D" = 165.58; "E" = 164.4 } }
    @{ Row = 64; Cells = @{ "A" = 45953; "D" = 161.47999999999999; "E" = 160.44999999999999; "F" = 170.45 } }
    @{ Row = 65; Cells = @{ "A" = 45953; "D" = 168.3; "E" = 171.03 } }
)

foreach ($update in $rowUpdates) {
    $r = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $ws.Range("$col$r").Value = $update.Cells[$col]
    }
}

